$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): reorder/rename columns and add new ones ---
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Nome"
$ws.Range("C1").Value = "CPF"
$ws.Range("D1").Value = "Data de Nascimento"
$ws.Range("E1").Value = "Senha"
$ws.Range("F1").Value = "Email"
$ws.Range("G1").Value = "autorizado?"
$ws.Range("H1").Value = "is_active"
$ws.Range("I1").Value = "is_authenticated"

# copy header style (bold/bordered/centered) onto the two newly added header cells
$ws.Range("A1").Copy()
$ws.Range("H1:I1").PasteSpecial(-4122)

# --- Row 2: existing "matheus" record, shifted into the new column layout ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "matheus"
$ws.Range("C2").Value = 48767507859
$ws.Range("D2").Value = 4041998
$ws.Range("E2").Value = 'pbkdf2:sha256:260000$FGVLB4z6awWJ2J7w$59be00fdf0984a3c0e0949d2b2846f8d0841e9d9e49de23dc2bfdfb39edda56c'
$ws.Range("F2").Value = "matheus@ufabc.edu.br"
$ws.Range("G2").Value = $true
$ws.Range("H2").Value = $true
$ws.Range("I2").Value = $true

# --- Row 3: new "teste" record ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "teste"
$ws.Range("C3").Value = 1234
$ws.Range("D3").Value = 1234
$ws.Range("E3").Value = 'pbkdf2:sha256:260000$5N0ofPnhoaIMtYtN$9a87a52797e5d0c6b0676b56e6e2888432a9940ebeee7a299824aac91ea32545'
$ws.Range("F3").Value = "teste@ufabc.com.br"
$ws.Range("G3").Value = $true
$ws.Range("H3").Value = $true
$ws.Range("I3").Value = $true
